# Adding a list of books + sql insert statements
# Builds a new column G with a CONCAT() formula that renders a SQL INSERT
# statement per book row, using the existing Title(A)/Author(B)/Price(E)/
# Date(F) columns, and stamps the header with a new "Database insert" label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header ---------------------------------------------------------------
$ws.Range("G1").Value = "Database insert"

# --- Row 2: first INSERT formula (kept as a standalone, non-shared formula,
# matching how Excel leaves the very first cell of a fill-down range) -------
$ws.Range("G2").Formula = '=CONCAT("INSERT INTO BOOK (TITLE, AUTHOR, PRICE, DATE) VALUES (''",A2,"'',''",B2,"'',''",E2,"'',''",F2,"''",")")'

# --- Rows 3-66: fill the formula down as one block (creates a shared
# formula group, same as dragging the fill handle down in the UI) ----------
$ws.Range("G3:G66").Formula = '=CONCAT("INSERT INTO BOOK (TITLE, AUTHOR, PRICE, DATE) VALUES (''",A3,"'',''",B3,"'',''",E3,"'',''",F3,"''",")")'

# --- Rows 67-101: second fill-down block (matches a later, separate
# extension of the range down to the last row) ------------------------------
$ws.Range("G67:G101").Formula = '=CONCAT("INSERT INTO BOOK (TITLE, AUTHOR, PRICE, DATE) VALUES (''",A67,"'',''",B67,"'',''",E67,"'',''",F67,"''",")")'

# --- Stray formatting artifact: a JetBrains Mono / blue / vertically
# centered font got tried out on a few empty cells past the data (L4:L6) --
$fmtRange = $ws.Range("L4:L6")
$fmtRange.Font.Name = "JetBrains Mono"
$fmtRange.Font.Size = 10
$fmtRange.Font.Color = 11744000
$fmtRange.VerticalAlignment = -4108

# --- Final selection left on the sheet -------------------------------------
$ws.Range("G2:G101").Select()
